$d = $word.ActiveDocument

# Step 1: Replace the placeholder text with the new line.
$d.Content.Find.Execute("(Will develop when making the program)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "List of techniques needed:", 2)

# Step 2: Append a brand-new, empty paragraph (same run formatting: Tahoma,
# size 24) right after the paragraph we just edited, i.e. at the very end
# of the document content.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$endRange.InsertXML($newParaXml)
